# Applies the "added to meeting outcomes" edit:
#  1. "Research Flask and MongoDb" -> split run + spellcheck proofErr around "MongoDb"
#  2. "More details ... look like" -> split run + gramStart/gramEnd proofErr around "look"
#  3. "3 4pm" -> split into "3 4p" / "m"
#  4. New "Goals:"/"Have skeletal..."/"Outcomes:"/"Decided..."/"Created SRS..." paragraphs
#     appended at the end of the document body.

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function New-PkgXml([string]$bodyFragment) {
    return @"
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="$wNs"><w:body>
$bodyFragment
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
}

# --- 1. "Research Flask and MongoDb" ------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Research Flask and MongoDb", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Text = ""
    $frag = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Research Flask and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MongoDb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
    $r.InsertXML((New-PkgXml $frag))
}

# --- 2. "More details ... look like" (keep the trailing "?" run) -------
$r = $d.Content
$found = $r.Find.Execute("More details on what each user would be able to do and what it would look like", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Text = ""
    $frag = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">More details on what each user would be able to do and what it would </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>look</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> like</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p>'
    $r.InsertXML((New-PkgXml $frag))
}

# --- 3 & 4. Split "3 4pm" and append the new Goals/Outcomes paragraphs -
$r = $d.Content
$found = $r.Find.Execute("4/13 4pm", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Text = ""
    $frag = @'
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs></w:pPr><w:r><w:t>4/1</w:t></w:r><w:r><w:t>3 4p</w:t></w:r><w:r><w:t>m</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs></w:pPr><w:r><w:t>Goals:</w:t></w:r></w:p>
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">Have skeletal implementation of perhaps Django web </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>app</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs><w:ind w:left="720"/></w:pPr></w:p>
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs><w:ind w:left="720"/></w:pPr><w:r><w:t>Outcomes:</w:t></w:r></w:p>
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">Decided on using Django as a framework for making the web </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>app</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">Created SRS and SDS Documents </w:t></w:r></w:p>
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1935"/></w:tabs><w:ind w:left="720"/></w:pPr></w:p>
'@
    $r.InsertXML((New-PkgXml $frag))
}

# InsertXML on the document-final paragraph leaves the original paragraph
# mark behind as a fresh trailing empty paragraph; give it back the
# ind:left=720 direct formatting the real edit ended up with.
$last = $d.Paragraphs.Last
$last.Range.ParagraphFormat.LeftIndent = 36
